# db data insertion into excel
# Add a new worksheet "dbdata" after the existing sheets, populate it with
# the "database" records (as text values, matching how the source data was
# pulled), and make it the active sheet.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "dbdata"

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
}

Set-TextValue $ws.Range("A1") "rakesh"
Set-TextValue $ws.Range("B1") "rakesh@gmail.com"
Set-TextValue $ws.Range("C1") "20000"

Set-TextValue $ws.Range("A2") "pavan"
Set-TextValue $ws.Range("B2") "pavan@gmail.com"
Set-TextValue $ws.Range("C2") "18000"

$ws.Columns.Item(1).ColumnWidth = 30.56
$ws.Columns.Item(2).ColumnWidth = 27.38
$ws.Columns.Item(3).ColumnWidth = 19.72

$ws.Range("A1").Select()
$ws.Select()
